$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 11).Value = 920.3333  # K11: 1078 -> 920.3333
$ws.Cells.Item(11, 13).Value = -780.3333  # M11: -938 -> -780.3333
$ws.Cells.Item(11, 8).Value = 920.3333  # H11: 1078 -> 920.3333
$ws.Cells.Item(11, 9).Value = 920.3333  # I11: 1078 -> 920.3333
$ws.Cells.Item(32, 8).Value = 500  # H32: 0 -> 500
$ws.Cells.Item(32, 13).Value = -174  # M32: None -> -174
$ws.Cells.Item(32, 9).Value = 500  # I32: 0 -> 500
$ws.Cells.Item(32, 11).Value = 500  # K32: 0 -> 500
$ws.Cells.Item(70, 13).Value = -5728.799999999999  # M70: -5578.5 -> -5728.799999999999
$ws.Cells.Item(70, 11).Value = 5998.799999999999  # K70: 5848.5 -> 5998.799999999999
$ws.Cells.Item(70, 9).Value = 1999.6  # I70: 1949.5 -> 1999.6
$ws.Cells.Item(70, 8).Value = 1999.6  # H70: 1949.5 -> 1999.6
$ws.Cells.Item(73, 11).Value = 5998.799999999999  # K73: 5848.5 -> 5998.799999999999
$ws.Cells.Item(73, 8).Value = 1999.6  # H73: 1949.5 -> 1999.6
$ws.Cells.Item(73, 13).Value = -5062.799999999999  # M73: -4912.5 -> -5062.799999999999
$ws.Cells.Item(73, 9).Value = 1999.6  # I73: 1949.5 -> 1999.6
$ws.Cells.Item(75, 12).Value = 0  # L75: 45000 -> 0
$ws.Cells.Item(75, 8).Value = 0  # H75: 45000 -> 0
$ws.Cells.Item(75, 14).Value = $null  # N75: remove (was -46872)
$ws.Cells.Item(75, 10).Value = 0  # J75: 45000 -> 0
$ws.Cells.Item(78, 10).Value = 0  # J78: 45000 -> 0
$ws.Cells.Item(78, 12).Value = 0  # L78: 135000 -> 0
$ws.Cells.Item(78, 14).Value = $null  # N78: remove (was -144360)
$ws.Cells.Item(78, 8).Value = 0  # H78: 45000 -> 0
$ws.Cells.Item(86, 13).Value = -2064  # M86: -1601.1667 -> -2064
$ws.Cells.Item(86, 11).Value = 3187  # K86: 2724.1667 -> 3187
$ws.Cells.Item(86, 9).Value = 3187  # I86: 2724.1667 -> 3187
$ws.Cells.Item(86, 8).Value = 6591.1665  # H86: 5393 -> 6591.1665
$ws.Cells.Item(89, 11).Value = 15935  # K89: 13620.8335 -> 15935
$ws.Cells.Item(89, 9).Value = 3187  # I89: 2724.1667 -> 3187
$ws.Cells.Item(89, 13).Value = -10319  # M89: -8004.833500000001 -> -10319
$ws.Cells.Item(89, 8).Value = 6591.1665  # H89: 5393 -> 6591.1665
$ws.Cells.Item(135, 11).Value = 7427.9997  # K135: 7483.5 -> 7427.9997
$ws.Cells.Item(135, 9).Value = 825.3333  # I135: 831.5 -> 825.3333
$ws.Cells.Item(135, 8).Value = 15496.429  # H135: 13665.625 -> 15496.429
$ws.Cells.Item(135, 13).Value = -4892.9997  # M135: -4948.5 -> -4892.9997
$ws.Cells.Item(137, 11).Value = 6067.200000000001  # K137: 5922 -> 6067.200000000001
$ws.Cells.Item(137, 14).Value = -25989.75  # N137: -24683.307 -> -25989.75
$ws.Cells.Item(137, 9).Value = 2022.4  # I137: 1974 -> 2022.4
$ws.Cells.Item(137, 8).Value = 4218.3335  # H137: 4015.3447 -> 4218.3335
$ws.Cells.Item(137, 10).Value = 6963.25  # J137: 6527.769 -> 6963.25
$ws.Cells.Item(137, 12).Value = 20889.75  # L137: 19583.307 -> 20889.75
$ws.Cells.Item(137, 13).Value = -3517.200000000001  # M137: -3372 -> -3517.200000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 14).Value = -4525  # N63: -4678 -> -4525
$ws.Cells.Item(63, 8).Value = 3065.889  # H63: 2826.3333 -> 3065.889
$ws.Cells.Item(63, 12).Value = 3153  # L63: 3306 -> 3153
$ws.Cells.Item(63, 13).Value = -2355  # M63: -2080.375 -> -2355
$ws.Cells.Item(63, 11).Value = 3041  # K63: 2766.375 -> 3041
$ws.Cells.Item(63, 10).Value = 3153  # J63: 3306 -> 3153
$ws.Cells.Item(63, 9).Value = 3041  # I63: 2766.375 -> 3041
$ws.Cells.Item(66, 12).Value = 15765  # L66: 16530 -> 15765
$ws.Cells.Item(66, 14).Value = -22629  # N66: -23394 -> -22629
$ws.Cells.Item(66, 11).Value = 15205  # K66: 13831.875 -> 15205
$ws.Cells.Item(66, 10).Value = 3153  # J66: 3306 -> 3153
$ws.Cells.Item(66, 8).Value = 3065.889  # H66: 2826.3333 -> 3065.889
$ws.Cells.Item(66, 9).Value = 3041  # I66: 2766.375 -> 3041
$ws.Cells.Item(66, 13).Value = -11773  # M66: -10399.875 -> -11773
$ws.Cells.Item(74, 10).Value = 19433.54  # J74: 19441.23 -> 19433.54
$ws.Cells.Item(74, 8).Value = 14721499  # H74: 14721505 -> 14721499
$ws.Cells.Item(74, 14).Value = -21181.54  # N74: -21189.23 -> -21181.54
$ws.Cells.Item(74, 12).Value = 19433.54  # L74: 19441.23 -> 19433.54
$ws.Cells.Item(77, 10).Value = 19433.54  # J77: 19441.23 -> 19433.54
$ws.Cells.Item(77, 14).Value = -105903.7  # N77: -105942.15 -> -105903.7
$ws.Cells.Item(77, 8).Value = 14721499  # H77: 14721505 -> 14721499
$ws.Cells.Item(77, 12).Value = 97167.70000000001  # L77: 97206.14999999999 -> 97167.70000000001
$ws.Cells.Item(97, 10).Value = 295.33334  # J97: 297.5 -> 295.33334
$ws.Cells.Item(97, 11).Value = 1108.1578  # K97: 1159.1666 -> 1108.1578
$ws.Cells.Item(97, 8).Value = 997.3182  # H97: 1073 -> 997.3182
$ws.Cells.Item(97, 12).Value = 295.33334  # L97: 297.5 -> 295.33334
$ws.Cells.Item(97, 13).Value = -612.1578  # M97: -663.1666 -> -612.1578
$ws.Cells.Item(97, 14).Value = -1287.33334  # N97: -1289.5 -> -1287.33334
$ws.Cells.Item(97, 9).Value = 1108.1578  # I97: 1159.1666 -> 1108.1578
$ws.Cells.Item(122, 11).Value = 4749.4998  # K122: 5474.25 -> 4749.4998
$ws.Cells.Item(122, 8).Value = 2381.125  # H122: 2564.1428 -> 2381.125
$ws.Cells.Item(122, 13).Value = -2299.4998  # M122: -3024.25 -> -2299.4998
$ws.Cells.Item(122, 9).Value = 1583.1666  # I122: 1824.75 -> 1583.1666
$ws.Cells.Item(132, 9).Value = 9998.75  # I132: 9999.333000000001 -> 9998.75
$ws.Cells.Item(132, 14).Value = -53863.10000000001  # N132: -55953.00199999999 -> -53863.10000000001
$ws.Cells.Item(132, 12).Value = 48803.10000000001  # L132: 50893.00199999999 -> 48803.10000000001
$ws.Cells.Item(132, 10).Value = 16267.7  # J132: 16964.334 -> 16267.7
$ws.Cells.Item(132, 11).Value = 29996.25  # K132: 29997.999 -> 29996.25
$ws.Cells.Item(132, 8).Value = 14476.571  # H132: 15223.083 -> 14476.571
$ws.Cells.Item(132, 13).Value = -27466.25  # M132: -27467.999 -> -27466.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 13).Value = -1019  # M3: -1051.091 -> -1019
$ws.Cells.Item(3, 9).Value = 1133  # I3: 1165.091 -> 1133
$ws.Cells.Item(3, 11).Value = 1133  # K3: 1165.091 -> 1133
$ws.Cells.Item(3, 8).Value = 1133  # H3: 1165.091 -> 1133
$ws.Cells.Item(20, 12).Value = 2651.5  # L20: 3419.8 -> 2651.5
$ws.Cells.Item(20, 8).Value = 3138.1765  # H20: 3858.6428 -> 3138.1765
$ws.Cells.Item(20, 10).Value = 2651.5  # J20: 3419.8 -> 2651.5
$ws.Cells.Item(20, 11).Value = 3403.6365  # K20: 4102.4443 -> 3403.6365
$ws.Cells.Item(20, 9).Value = 3403.6365  # I20: 4102.4443 -> 3403.6365
$ws.Cells.Item(20, 13).Value = -3156.6365  # M20: -3855.4443 -> -3156.6365
$ws.Cells.Item(20, 14).Value = -3145.5  # N20: -3913.8 -> -3145.5
$ws.Cells.Item(86, 13).Value = -1390.5  # M86: -3603.5 -> -1390.5
$ws.Cells.Item(86, 11).Value = 2513.5  # K86: 4726.5 -> 2513.5
$ws.Cells.Item(86, 12).Value = 3700  # L86: 0 -> 3700
$ws.Cells.Item(86, 9).Value = 2513.5  # I86: 4726.5 -> 2513.5
$ws.Cells.Item(86, 10).Value = 3700  # J86: 0 -> 3700
$ws.Cells.Item(86, 14).Value = -5946  # N86: None -> -5946
$ws.Cells.Item(86, 8).Value = 2645.3333  # H86: 4726.5 -> 2645.3333
$ws.Cells.Item(89, 12).Value = 18500  # L89: 0 -> 18500
$ws.Cells.Item(89, 11).Value = 12567.5  # K89: 23632.5 -> 12567.5
$ws.Cells.Item(89, 10).Value = 3700  # J89: 0 -> 3700
$ws.Cells.Item(89, 9).Value = 2513.5  # I89: 4726.5 -> 2513.5
$ws.Cells.Item(89, 14).Value = -29732  # N89: None -> -29732
$ws.Cells.Item(89, 13).Value = -6951.5  # M89: -18016.5 -> -6951.5
$ws.Cells.Item(89, 8).Value = 2645.3333  # H89: 4726.5 -> 2645.3333
$ws.Cells.Item(92, 12).Value = 64156  # L92: 68694.25 -> 64156
$ws.Cells.Item(92, 10).Value = 64156  # J92: 68694.25 -> 64156
$ws.Cells.Item(92, 8).Value = 64156  # H92: 68694.25 -> 64156
$ws.Cells.Item(92, 14).Value = -69148  # N92: -73686.25 -> -69148
$ws.Cells.Item(105, 10).Value = 2531.6667  # J105: 3500 -> 2531.6667
$ws.Cells.Item(105, 8).Value = 1978.8667  # H105: 2855.4443 -> 1978.8667
$ws.Cells.Item(105, 14).Value = -6025.6667  # N105: -6994 -> -6025.6667
$ws.Cells.Item(105, 13).Value = 597.3334  # M105: 180.6666 -> 597.3334
$ws.Cells.Item(105, 12).Value = 2531.6667  # L105: 3500 -> 2531.6667
$ws.Cells.Item(105, 11).Value = 1149.6666  # K105: 1566.3334 -> 1149.6666
$ws.Cells.Item(105, 9).Value = 1149.6666  # I105: 1566.3334 -> 1149.6666
$ws.Cells.Item(107, 14).Value = -5322  # N107: None -> -5322
$ws.Cells.Item(107, 12).Value = 1482  # L107: 0 -> 1482
$ws.Cells.Item(107, 13).Value = 398.1111000000001  # M107: 397.3334 -> 398.1111000000001
$ws.Cells.Item(107, 9).Value = 1521.8889  # I107: 1522.6666 -> 1521.8889
$ws.Cells.Item(107, 11).Value = 1521.8889  # K107: 1522.6666 -> 1521.8889
$ws.Cells.Item(107, 10).Value = 1482  # J107: 0 -> 1482
$ws.Cells.Item(107, 8).Value = 1517.9  # H107: 1522.6666 -> 1517.9
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 9).Value = 811  # I58: 764.1667 -> 811
$ws.Cells.Item(58, 11).Value = 811  # K58: 764.1667 -> 811
$ws.Cells.Item(58, 8).Value = 1712  # H58: 1180.2222 -> 1712
$ws.Cells.Item(58, 13).Value = -608  # M58: -561.1667 -> -608
$ws.Cells.Item(118, 14).Value = -91303  # N118: -98303.664 -> -91303
$ws.Cells.Item(118, 12).Value = 87989  # L118: 94989.664 -> 87989
$ws.Cells.Item(118, 8).Value = 87989  # H118: 94989.664 -> 87989
$ws.Cells.Item(118, 10).Value = 87989  # J118: 94989.664 -> 87989
$ws.Cells.Item(132, 9).Value = 1895.6428  # I132: 1979.9231 -> 1895.6428
$ws.Cells.Item(132, 11).Value = 5686.928400000001  # K132: 5939.7693 -> 5686.928400000001
$ws.Cells.Item(132, 8).Value = 1832.4  # H132: 1906.1428 -> 1832.4
$ws.Cells.Item(132, 13).Value = -3156.928400000001  # M132: -3409.7693 -> -3156.928400000001
$ws.Cells.Item(134, 13).Value = -3343764  # M134: -3761117.7 -> -3343764
$ws.Cells.Item(134, 9).Value = 1115433  # I134: 1254550.9 -> 1115433
$ws.Cells.Item(134, 8).Value = 722994.9  # H134: 778418.3 -> 722994.9
$ws.Cells.Item(134, 11).Value = 3346299  # K134: 3763652.7 -> 3346299
$ws.Cells.Item(136, 8).Value = 1712  # H136: 1180.2222 -> 1712
$ws.Cells.Item(136, 13).Value = 117  # M136: 257.4998999999998 -> 117
$ws.Cells.Item(136, 9).Value = 811  # I136: 764.1667 -> 811
$ws.Cells.Item(136, 11).Value = 2433  # K136: 2292.5001 -> 2433
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 9).Value = 5073.4287  # I129: 935.6667 -> 5073.4287
$ws.Cells.Item(129, 12).Value = 250722630  # L129: 334293810 -> 250722630
$ws.Cells.Item(129, 11).Value = 15220.2861  # K129: 2807.0001 -> 15220.2861
$ws.Cells.Item(129, 14).Value = -250732630  # N129: -334303810 -> -250732630
$ws.Cells.Item(129, 10).Value = 83574210  # J129: 111431270 -> 83574210
$ws.Cells.Item(129, 13).Value = -10220.2861  # M129: 2192.9999 -> -10220.2861
$ws.Cells.Item(129, 8).Value = 30393850  # H129: 37144380 -> 30393850
$ws.Cells.Item(134, 14).Value = -31740  # N134: None -> -31740
$ws.Cells.Item(134, 10).Value = 7200  # J134: 0 -> 7200
$ws.Cells.Item(134, 13).Value = -3049.5651  # M134: -4044.428400000001 -> -3049.5651
$ws.Cells.Item(134, 12).Value = 21600  # L134: 0 -> 21600
$ws.Cells.Item(134, 9).Value = 2706.5217  # I134: 3038.1428 -> 2706.5217
$ws.Cells.Item(134, 8).Value = 3508.9285  # H134: 3038.1428 -> 3508.9285
$ws.Cells.Item(134, 11).Value = 8119.5651  # K134: 9114.428400000001 -> 8119.5651
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(55, 11).Value = 5500  # K55: 6000 -> 5500
$ws.Cells.Item(55, 13).Value = -5173  # M55: -5673 -> -5173
$ws.Cells.Item(55, 8).Value = 7000  # H55: 8000 -> 7000
$ws.Cells.Item(55, 9).Value = 5500  # I55: 6000 -> 5500
$ws.Cells.Item(80, 13).Value = -10772  # M80: -3002 -> -10772
$ws.Cells.Item(80, 11).Value = 11770  # K80: 4000 -> 11770
$ws.Cells.Item(80, 9).Value = 11770  # I80: 4000 -> 11770
$ws.Cells.Item(80, 8).Value = 11770  # H80: 4000 -> 11770
$ws.Cells.Item(83, 8).Value = 11770  # H83: 4000 -> 11770
$ws.Cells.Item(83, 11).Value = 58850  # K83: 20000 -> 58850
$ws.Cells.Item(83, 13).Value = -53858  # M83: -15008 -> -53858
$ws.Cells.Item(83, 9).Value = 11770  # I83: 4000 -> 11770
$ws.Cells.Item(97, 10).Value = 912.75  # J97: 723.3333 -> 912.75
$ws.Cells.Item(97, 11).Value = 1106.1177  # K97: 1664.4736 -> 1106.1177
$ws.Cells.Item(97, 8).Value = 1069.2858  # H97: 1438.6 -> 1069.2858
$ws.Cells.Item(97, 12).Value = 912.75  # L97: 723.3333 -> 912.75
$ws.Cells.Item(97, 13).Value = -610.1177  # M97: -1168.4736 -> -610.1177
$ws.Cells.Item(97, 14).Value = -1904.75  # N97: -1715.3333 -> -1904.75
$ws.Cells.Item(97, 9).Value = 1106.1177  # I97: 1664.4736 -> 1106.1177
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 13).Value = -770.6  # M61: -851.6666 -> -770.6
$ws.Cells.Item(61, 8).Value = 978.0714  # H61: 1028.9 -> 978.0714
$ws.Cells.Item(61, 9).Value = 972.6  # I61: 1053.6666 -> 972.6
$ws.Cells.Item(61, 11).Value = 972.6  # K61: 1053.6666 -> 972.6
$ws.Cells.Item(82, 11).Value = 516  # K82: 510.25 -> 516
$ws.Cells.Item(82, 9).Value = 516  # I82: 510.25 -> 516
$ws.Cells.Item(82, 12).Value = 1800  # L82: 1267.6666 -> 1800
$ws.Cells.Item(82, 14).Value = -2522  # N82: -1989.6666 -> -2522
$ws.Cells.Item(82, 13).Value = -155  # M82: -149.25 -> -155
$ws.Cells.Item(82, 10).Value = 1800  # J82: 1267.6666 -> 1800
$ws.Cells.Item(82, 8).Value = 1029.6  # H82: 834.8570999999999 -> 1029.6
$ws.Cells.Item(85, 14).Value = -4296  # N85: -3763.6666 -> -4296
$ws.Cells.Item(85, 13).Value = 732  # M85: 737.75 -> 732
$ws.Cells.Item(85, 10).Value = 1800  # J85: 1267.6666 -> 1800
$ws.Cells.Item(85, 8).Value = 1029.6  # H85: 834.8570999999999 -> 1029.6
$ws.Cells.Item(85, 11).Value = 516  # K85: 510.25 -> 516
$ws.Cells.Item(85, 12).Value = 1800  # L85: 1267.6666 -> 1800
$ws.Cells.Item(85, 9).Value = 516  # I85: 510.25 -> 516
$ws.Cells.Item(93, 8).Value = 142860050  # H93: 166669660 -> 142860050
$ws.Cells.Item(93, 13).Value = -166668062  # M93: -200001442 -> -166668062
$ws.Cells.Item(93, 11).Value = 166669310  # K93: 200002690 -> 166669310
$ws.Cells.Item(93, 9).Value = 166669310  # I93: 200002690 -> 166669310
$ws.Cells.Item(98, 14).Value = -118990  # N98: -119990 -> -118990
$ws.Cells.Item(98, 12).Value = 113000  # L98: 114000 -> 113000
$ws.Cells.Item(98, 8).Value = 113000  # H98: 114000 -> 113000
$ws.Cells.Item(98, 10).Value = 113000  # J98: 114000 -> 113000
$ws.Cells.Item(113, 13).Value = 1197.4  # M113: 1116.3334 -> 1197.4
$ws.Cells.Item(113, 9).Value = 972.6  # I113: 1053.6666 -> 972.6
$ws.Cells.Item(113, 11).Value = 972.6  # K113: 1053.6666 -> 972.6
$ws.Cells.Item(113, 8).Value = 978.0714  # H113: 1028.9 -> 978.0714
$ws.Cells.Item(132, 9).Value = 7063.5264  # I132: 7101.6313 -> 7063.5264
$ws.Cells.Item(132, 11).Value = 21190.5792  # K132: 21304.8939 -> 21190.5792
$ws.Cells.Item(132, 8).Value = 58486.91  # H132: 58519.816 -> 58486.91
$ws.Cells.Item(132, 13).Value = -18660.5792  # M132: -18774.8939 -> -18660.5792
$ws.Cells.Item(136, 8).Value = 33044.227  # H136: 33010.45 -> 33044.227
$ws.Cells.Item(136, 10).Value = 80558.8  # J136: 86152.28999999999 -> 80558.8
$ws.Cells.Item(136, 13).Value = -11056.44  # M136: -10636.845 -> -11056.44
$ws.Cells.Item(136, 9).Value = 4535.48  # I136: 4395.615 -> 4535.48
$ws.Cells.Item(136, 14).Value = -246776.4  # N136: -263556.87 -> -246776.4
$ws.Cells.Item(136, 12).Value = 241676.4  # L136: 258456.87 -> 241676.4
$ws.Cells.Item(136, 11).Value = 13606.44  # K136: 13186.845 -> 13606.44
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 14).Value = -17355.666  # N41: -18460.5 -> -17355.666
$ws.Cells.Item(41, 10).Value = 16575.666  # J41: 17680.5 -> 16575.666
$ws.Cells.Item(41, 8).Value = 16575.666  # H41: 17680.5 -> 16575.666
$ws.Cells.Item(41, 12).Value = 16575.666  # L41: 17680.5 -> 16575.666
$ws.Cells.Item(62, 9).Value = 7743.2163  # I62: 7900 -> 7743.2163
$ws.Cells.Item(62, 11).Value = 7743.2163  # K62: 7900 -> 7743.2163
$ws.Cells.Item(62, 14).Value = -28576634  # N62: -40005808 -> -28576634
$ws.Cells.Item(62, 13).Value = -7119.2163  # M62: -7276 -> -7119.2163
$ws.Cells.Item(62, 10).Value = 28575386  # J62: 40004560 -> 28575386
$ws.Cells.Item(62, 8).Value = 4552595.5  # H62: 5007482.5 -> 4552595.5
$ws.Cells.Item(62, 12).Value = 28575386  # L62: 40004560 -> 28575386
$ws.Cells.Item(65, 10).Value = 28575386  # J65: 40004560 -> 28575386
$ws.Cells.Item(65, 13).Value = -35596.0815  # M65: -36380 -> -35596.0815
$ws.Cells.Item(65, 9).Value = 7743.2163  # I65: 7900 -> 7743.2163
$ws.Cells.Item(65, 12).Value = 142876930  # L65: 200022800 -> 142876930
$ws.Cells.Item(65, 8).Value = 4552595.5  # H65: 5007482.5 -> 4552595.5
$ws.Cells.Item(65, 14).Value = -142883170  # N65: -200029040 -> -142883170
$ws.Cells.Item(65, 11).Value = 38716.0815  # K65: 39500 -> 38716.0815
$ws.Cells.Item(122, 11).Value = 6719.2002  # K122: 6985.071599999999 -> 6719.2002
$ws.Cells.Item(122, 8).Value = 5654.7666  # H122: 5815.3105 -> 5654.7666
$ws.Cells.Item(122, 13).Value = -4269.2002  # M122: -4535.071599999999 -> -4269.2002
$ws.Cells.Item(122, 9).Value = 2239.7334  # I122: 2328.3572 -> 2239.7334
